$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stipend text in cell B5 with the new amounts (3200 / 5600 rubles
# instead of the previous 3383 / 5751 rubles).
$newText = "Стипендия есть у студентов, учащихся на бюджете, до первой сессии она есть у всех первокурсников и равна стипендии для студентов, сдавших сессию на 4 и 5. После первой сессии для студентов, сдавших сессию на 4 и 5 она составляет, по последним данным, 3200 рублей, для сдавших сессию на 5 - 5600 рублей. Сумма стипендии иногда меняется, поэтому советую отслеживать новости в группе профкома фсир в вк"

$ws.Range("B5").Value = $newText

# Move the viewport/selection like the author did while editing this cell.
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 2
